# Adjusted risk calc formula
# For the three "LF" rows below, the Total Risk column changes from 4 to 6
# and the Current Risk column changes from L to M (Future Risk stays L).
#   - LF1: predation from pinnipeds or other aquatic species
#   - LF41: competition with hatchery fry
#   - LF69: rearing in a hatchery environment leading to maladaptation

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$targets = @(
    "LF1: Mortality or fitness reduction due to predation from pinnipeds or other aquatic species",
    "LF41: Mortality or fitness reduction as a result of competition with hatchery fry",
    "LF69: Mortality or fitness reduction as a result of rearing in a hatchery environment leading to maladaptation to the wild environment. This is measured in a reduction in PNI."
)

$rowCount = $table.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $row = $table.Rows.Item($r)
    $lfText = $row.Cells.Item(2).Range.Text
    $lfText = $lfText.TrimEnd([char]13, [char]7).Trim()

    foreach ($target in $targets) {
        if ($lfText -eq $target) {
            # Column 4 = Total Risk, Column 5 = Current Risk
            $totalRiskCell = $row.Cells.Item(4)
            $totalRiskCell.Range.Text = "6"

            $currentRiskCell = $row.Cells.Item(5)
            $currentRiskCell.Range.Text = "M"
        }
    }
}
